$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old "avg_commute" column (O) entirely.
# ---------------------------------------------------------------------------
$ws.Range("O1:O4").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. Column A header becomes "Commute" (keeps the bold/bordered header
#    style used by the rest of row 1, by copying B1's formatting over).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Commute"
$ws.Cells.Item(1, 2).Copy() | Out-Null
$ws.Cells.Item(1, 1).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Column A data cells (rows 2-4) switch from a numeric index to a literal
#    text "TRUE" flag, and they lose the bordered/bold header style that
#    used to be applied to them (plain default formatting, like column B).
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "'TRUE"
$ws.Cells.Item(2, 1).ClearFormats()
$ws.Cells.Item(3, 1).Value = "'TRUE"
$ws.Cells.Item(3, 1).ClearFormats()
$ws.Cells.Item(4, 1).Value = "'TRUE"
$ws.Cells.Item(4, 1).ClearFormats()

# ---------------------------------------------------------------------------
# 4. Update the numeric values for rows 2-4 (columns C..L changed; B, M, N
#    stay the same).
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 3).Value = 50460
$ws.Cells.Item(2, 4).Value = 0.74853024584738
$ws.Cells.Item(2, 5).Value = 0.3526294182996187
$ws.Cells.Item(2, 6).Value = 7.812576636820745
$ws.Cells.Item(2, 7).Value = 246.9
$ws.Cells.Item(2, 8).Value = 1518.5
$ws.Cells.Item(2, 9).Value = 3527.5
$ws.Cells.Item(2, 10).Value = 3789.3
$ws.Cells.Item(2, 11).Value = 1794.3
$ws.Cells.Item(2, 12).Value = 39524.4226115161

$ws.Cells.Item(3, 3).Value = 47241
$ws.Cells.Item(3, 4).Value = 0.733787678404157
$ws.Cells.Item(3, 5).Value = 0.3823043262558196
$ws.Cells.Item(3, 6).Value = 7.569324569233842
$ws.Cells.Item(3, 7).Value = 303.2
$ws.Cells.Item(3, 8).Value = 1539.9
$ws.Cells.Item(3, 9).Value = 3184.2
$ws.Cells.Item(3, 10).Value = 3474.4
$ws.Cells.Item(3, 11).Value = 1807.8
$ws.Cells.Item(3, 12).Value = 35744.22831054959

$ws.Cells.Item(4, 3).Value = -6.379310344827586
$ws.Cells.Item(4, 4).Value = -1.969535302682866
$ws.Cells.Item(4, 5).Value = 8.415323967947291
$ws.Cells.Item(4, 6).Value = -3.113595922252512
$ws.Cells.Item(4, 7).Value = 22.80275415147833
$ws.Cells.Item(4, 8).Value = 1.409285479091215
$ws.Cells.Item(4, 9).Value = -9.732104890148836
$ws.Cells.Item(4, 10).Value = -8.310241997202651
$ws.Cells.Item(4, 11).Value = 0.7523825447249624
$ws.Cells.Item(4, 12).Value = -9.564198668053633

# ---------------------------------------------------------------------------
# 5. Three brand-new rows (5, 6, 7) holding the "FALSE" (commute = no) block,
#    mirroring the structure of rows 2-4, with plain/default formatting.
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 1).Value = "'FALSE"
$ws.Cells.Item(5, 2).Value = 2
$ws.Cells.Item(5, 3).Value = 199460
$ws.Cells.Item(5, 4).Value = 0.6422447850959879
$ws.Cells.Item(5, 5).Value = 0.3159372780047863
$ws.Cells.Item(5, 6).Value = 7.280392554873131
$ws.Cells.Item(5, 7).Value = 149.9622641509434
$ws.Cells.Item(5, 8).Value = 1324.528301886792
$ws.Cells.Item(5, 9).Value = 2438.867924528302
$ws.Cells.Item(5, 10).Value = 2415.075471698113
$ws.Cells.Item(5, 11).Value = 1191.471698113208
$ws.Cells.Item(5, 12).Value = 27360.44954639928
$ws.Cells.Item(5, 13).Value = 1
$ws.Cells.Item(5, 14).Value = 0

$ws.Cells.Item(6, 1).Value = "'FALSE"
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(6, 3).Value = 197204
$ws.Cells.Item(6, 4).Value = 0.6036108708879523
$ws.Cells.Item(6, 5).Value = 0.3426993618394082
$ws.Cells.Item(6, 6).Value = 6.877671503647454
$ws.Cells.Item(6, 7).Value = 168.7924528301887
$ws.Cells.Item(6, 8).Value = 1445.754716981132
$ws.Cells.Item(6, 9).Value = 2275.075471698113
$ws.Cells.Item(6, 10).Value = 2242.811320754717
$ws.Cells.Item(6, 11).Value = 1272.811320754717
$ws.Cells.Item(6, 12).Value = 25567.91448016014
$ws.Cells.Item(6, 13).Value = 1
$ws.Cells.Item(6, 14).Value = 1

$ws.Cells.Item(7, 1).Value = "'FALSE"
$ws.Cells.Item(7, 2).Value = "Difference (%)"
$ws.Cells.Item(7, 3).Value = -1.131053845382533
$ws.Cells.Item(7, 4).Value = -6.01545004406093
$ws.Cells.Item(7, 5).Value = 8.4706951973602
$ws.Cells.Item(7, 6).Value = -5.53158429563137
$ws.Cells.Item(7, 7).Value = 12.55661801711122
$ws.Cells.Item(7, 8).Value = 9.152421652421658
$ws.Cells.Item(7, 9).Value = -6.715921398731241
$ws.Cells.Item(7, 10).Value = -7.132868225533007
$ws.Cells.Item(7, 11).Value = 6.826819535060483
$ws.Cells.Item(7, 12).Value = -6.5515556065673
$ws.Cells.Item(7, 13).Value = 0
$ws.Cells.Item(7, 14).Value = "inf"

# The new rows have no special styling (default style) - match that.
$ws.Range("A5:N7").ClearFormats()

$wb.Save()
